$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 9 (A9:D9) currently uses the "pending" (orange) style shared with the
# header-like rows; it needs to move to the "done" (green) style used by the
# surrounding data rows (e.g. row 10), and its Situation changes to
# "Desenvolvido" while the delivery forecast becomes a concrete date.
$ws.Range("A10:D10").Copy()
$ws.Range("A9:D9").PasteSpecial(-4122)

$ws.Range("C9").Value = "Desenvolvido"
$ws.Range("D9").Value = 41315

# Rows 10-14: previsao de entrega date changes from 07/12/2013 (41615) to
# 10/02/2013 (41315)
$ws.Range("D10").Value = 41315
$ws.Range("D11").Value = 41315
$ws.Range("D12").Value = 41315
$ws.Range("D13").Value = 41315
$ws.Range("D14").Value = 41315

# Update the active selection shown when the workbook is reopened.
$ws.Range("G16").Select()

Write-Host "done"
